$d = $word.ActiveDocument

# --- Edit 1 ---------------------------------------------------------------
# Merge the paragraph ending "...generar una notificacion " with the
# (empty) paragraph after it and the paragraph starting " que son una
# buena herramienta..." into a single paragraph, while rewriting the
# wording about alert systems / sentiment-analysis systems.
$rng1 = $d.Content
$rng1.Find.ClearFormatting()
$rng1.Find.Replacement.ClearFormatting()
$old1 = " mecanismos como los sistemas de alerta. Los sistemas de alerta, son sistemas que ante uno o varios eventos, son capaces de generar una notificación ^p^p que son una buena herramienta para el filtrado de la información, porque a través de un conjunto de procedimientos, son capaces de seleccionar la información y procesarla. Otros mecanismos, como los sistemas de análisis de sentimiento"
$new1 = " mecanismos como los sistemas de alerta, sistemas de análisis de sentimientos. Los sistemas de alerta, son sistemas que ante uno o varios eventos, son capaces de generar una notificación o alerta, para indicar a un usuario o grupo  de usuarios sobre un suceso, acontecimiento o hacer una sugerencia. Los sistemas de análisis de sentimiento"
$ok1 = $rng1.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Edit 2 ---------------------------------------------------------------
# Insert the clarification about "mineria de opinion" right after
# "...sentimientos" and before the following ", utilizan...".
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Replacement.ClearFormatting()
$old2 = "sentimientos, utilizan"
$new2 = "sentimientos o también conocidos como “minería de opinión”, utilizan"
$ok2 = $rng2.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- Edit 3 ---------------------------------------------------------------
# "Estas herramientas, aplicadas" -> "Estos sistemas, aplicados"
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Replacement.ClearFormatting()
$old3 = "Estas herramientas, aplicadas en el sector turístico"
$new3 = "Estos sistemas, aplicados en el sector turístico"
$ok3 = $rng3.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

if (-not $ok1) { throw "Edit 1 (merge alert/sentiment-system paragraphs) did not find its target text." }
if (-not $ok2) { throw "Edit 2 (mineria de opinion insertion) did not find its target text." }
if (-not $ok3) { throw "Edit 3 (Estas herramientas -> Estos sistemas) did not find its target text." }

Write-Output "Edit1: $ok1 Edit2: $ok2 Edit3: $ok3"
